$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.894.43"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +4.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.772.68"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +4.81%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.76"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +3.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.01"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +4.20%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +4.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.61"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +5.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.03"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.57"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.206.44"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +4.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.763.32"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +5.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "51.723.89"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.874"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.20"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +10.54%  "
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "275.09"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.72"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  +5.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.63"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.17"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.67"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0819"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.90"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.09"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.92"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0378"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +8.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +24.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.38"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.115"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.60"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.18%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.00"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.065.66"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.51"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +5.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.85"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.96"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.42%  "
